$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# Row 9 is a service-log entry that was filled in with real data.
# D9 ("1033") looks numeric, so force it to stay text like the rest of
# the sheet (every other cell in this table is stored as text).
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1033"
$ws.Range("D9").Style = "Normal"

$ws.Range("L9").Value = "26/2/2026"
$ws.Range("M9").Value = "سيرفيس"
$ws.Range("N9").Value = "تم تغير سلك سلندر"
$ws.Range("O9").Value = "م.سامر.م.محمد عبدالله.محمود إيهاب.حسام .سعيد .محمد ابراهيم "
$ws.Range("Q9").Value = "Done "

# Column Q ("Cylinder (o)") rows 2-34 get filled in with "nan" placeholders
# (row 9's Q cell is handled above with the real "Done " value).
foreach ($r in 2..34) {
    if ($r -ne 9) {
        $ws.Range("Q$r").Value = "nan"
    }
}
